$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings (e.g. "1.0000", "26.842.93")
# are preserved exactly as text and are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.842.93"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").Value = "1.873.48"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").Value = "301.03"
$ws.Range("E5").Value = "  -2.13%  "

$ws.Range("D7").Value = "0.5325"
$ws.Range("E7").Value = "  +1.26%  "

$ws.Range("D8").Value = "0.3753"
$ws.Range("E8").Value = "  -1.43%  "

$ws.Range("D9").Value = "0.07178"
$ws.Range("E9").Value = "  -1.78%  "

$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("D11").Value = "0.8880"
$ws.Range("E11").Value = "  -1.90%  "

$ws.Range("D12").Value = "0.08156"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").Value = "1.912.01"
$ws.Range("E13").Value = "  +4.32%  "

$ws.Range("E14").Value = "  -2.54%  "

$ws.Range("D15").Value = "5.302"
$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("D18").Value = "0.000008555"
$ws.Range("E18").Value = "  -1.46%  "

$ws.Range("D20").Value = "26.870.93"
$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("D21").Value = "4.988"
$ws.Range("E21").Value = "  -2.73%  "

$ws.Range("D22").Value = "10.67"
$ws.Range("E22").Value = "  -1.42%  "

$ws.Range("D23").Value = "6.395"
$ws.Range("E23").Value = "  -1.34%  "

$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "146.34"
$ws.Range("E24").Value = "  -2.16%  "

$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "2.278"
$ws.Range("E25").Value = "  -3.44%  "

$ws.Range("D26").Value = "1.732"
$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("D28").Value = "113.70"
$ws.Range("E28").Value = "  -3.03%  "

$ws.Range("D29").Value = "4.729"
$ws.Range("E29").Value = "  -2.41%  "

$ws.Range("D30").Value = "4.614"
$ws.Range("E30").Value = "  -5.64%  "

$ws.Range("D31").Value = "0.09133"
$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("D32").Value = "0.8157"
$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("D33").Value = "0.04991"
$ws.Range("E33").Value = "  -1.56%  "

$ws.Range("D34").Value = "1.176"
$ws.Range("E34").Value = "  -4.36%  "

$ws.Range("E35").Value = "  -0.66%  "

$ws.Range("D36").Value = "0.6082"
$ws.Range("E36").Value = "  +6.28%  "

$ws.Range("E37").Value = "  -5.97%  "

$ws.Range("D38").Value = "2.598"
$ws.Range("E38").Value = "  -3.63%  "

$ws.Range("D39").Value = "0.01951"
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.590"
$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "8.903"
$ws.Range("E42").Value = "  -1.12%  "

$ws.Range("E43").Value = "  +4.87%  "

$ws.Range("D44").Value = "114.87"
$ws.Range("E44").Value = "  -1.50%  "

$ws.Range("E45").Value = "  -1.53%  "

$ws.Range("D46").Value = "1.0000"
$ws.Range("E46").Value = "  -0.26%  "

$ws.Range("D47").Value = "1.631"
$ws.Range("E47").Value = "  -0.67%  "

$ws.Range("D48").Value = "9.896"
$ws.Range("E48").Value = "  -2.43%  "

$ws.Range("D49").Value = "37.51"
$ws.Range("E49").Value = "  -2.73%  "

$ws.Range("D50").Value = "0.06064"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("D51").Value = "62.20"
$ws.Range("E51").Value = "  -3.27%  "
